$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.35876433333334
$ws.Range("H2").Value = 58.07629300000001
$ws.Range("I2").Value = 0.9707667559429034
$ws.Range("J2").Value = 0.9707667559429034
$ws.Range("M2").Value = 7.413580666666667
$ws.Range("N2").Value = 22.240742
$ws.Range("O2").Value = 0.05108888817597561
$ws.Range("P2").Value = 0.05108888817597561
$ws.Range("Q2").Value = 143.5177609921563
$ws.Range("R2").Value = 1291.659848929406
$ws.Range("S2").Value = 0.04959539423932161
$ws.Range("T2").Value = 0.0495953942393216
$ws.Range("G3").Value = 19.35876433333334
$ws.Range("H3").Value = 58.07629300000001
$ws.Range("I3").Value = 0.9707667559429034
$ws.Range("J3").Value = 0.9707667559429034
$ws.Range("O3").Value = 0.0112127179963522
$ws.Range("P3").Value = 0.0112127179963522
$ws.Range("Q3").Value = 31.49851638833779
$ws.Range("R3").Value = 283.4866474950401
$ws.Range("S3").Value = 0.01088493387462144
$ws.Range("T3").Value = 0.01088493387462144
$ws.Range("G4").Value = 19.35876433333334
$ws.Range("H4").Value = 58.07629300000001
$ws.Range("I4").Value = 0.9707667559429034
$ws.Range("J4").Value = 0.9707667559429034
$ws.Range("O4").Value = 0.9376983938276722
$ws.Range("P4").Value = 0.9376983938276722
$ws.Range("Q4").Value = 2634.161336698902
$ws.Range("R4").Value = 23707.45203029012
$ws.Range("S4").Value = 0.9102864278289603
$ws.Range("T4").Value = 0.9102864278289603
$ws.Range("H5").Value = 0.9049070000000001
$ws.Range("I5").Value = 0.01512585579145048
$ws.Range("J5").Value = 0.01512585579145048
$ws.Range("M5").Value = 7.413580666666667
$ws.Range("N5").Value = 22.240742
$ws.Range("O5").Value = 0.05108888817597561
$ws.Range("P5").Value = 0.05108888817597561
$ws.Range("Q5").Value = 2.236200346777112
$ws.Range("R5").Value = 20.125803120994
$ws.Range("S5").Value = 0.0007727631550953465
$ws.Range("T5").Value = 0.0007727631550953465
$ws.Range("H6").Value = 0.9049070000000001
$ws.Range("I6").Value = 0.01512585579145048
$ws.Range("J6").Value = 0.01512585579145048
$ws.Range("O6").Value = 0.0112127179963522
$ws.Range("P6").Value = 0.0112127179963522
$ws.Range("R6").Value = 4.417104440960001
$ws.Range("S6").Value = 0.0001696019554430249
$ws.Range("T6").Value = 0.000169601955443025
$ws.Range("H7").Value = 0.9049070000000001
$ws.Range("I7").Value = 0.01512585579145048
$ws.Range("J7").Value = 0.01512585579145048
$ws.Range("O7").Value = 0.9376983938276722
$ws.Range("P7").Value = 0.9376983938276722
$ws.Range("S7").Value = 0.0141834906809121
$ws.Range("T7").Value = 0.01418349068091211
$ws.Range("G8").Value = 0.2813256666666666
$ws.Range("H8").Value = 0.843977
$ws.Range("I8").Value = 0.01410738826564608
$ws.Range("J8").Value = 0.01410738826564608
$ws.Range("M8").Value = 7.413580666666667
$ws.Range("N8").Value = 22.240742
$ws.Range("O8").Value = 0.05108888817597561
$ws.Range("P8").Value = 0.05108888817597561
$ws.Range("Q8").Value = 2.085630523437111
$ws.Range("R8").Value = 18.770674710934
$ws.Range("S8").Value = 0.0007207307815586631
$ws.Range("T8").Value = 0.0007207307815586631
$ws.Range("G9").Value = 0.2813256666666666
$ws.Range("H9").Value = 0.843977
$ws.Range("I9").Value = 0.01410738826564608
$ws.Range("J9").Value = 0.01410738826564608
$ws.Range("O9").Value = 0.0112127179963522
$ws.Range("P9").Value = 0.0112127179963522
$ws.Range("Q9").Value = 0.4577431167288889
$ws.Range("R9").Value = 4.11968805056
$ws.Range("S9").Value = 0.0001581821662877377
$ws.Range("T9").Value = 0.0001581821662877377
$ws.Range("G10").Value = 0.2813256666666666
$ws.Range("H10").Value = 0.843977
$ws.Range("I10").Value = 0.01410738826564608
$ws.Range("J10").Value = 0.01410738826564608
$ws.Range("O10").Value = 0.9376983938276722
$ws.Range("P10").Value = 0.9376983938276722
$ws.Range("R10").Value = 344.521718047124
$ws.Range("S10").Value = 0.01322847531779968
$ws.Range("T10").Value = 0.01322847531779968
